$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A315").Value = "2023-12-14 12:01:54"
$ws.Range("B315").Value = 0.0008

$ws.Range("A316").Value = "2023-12-14 12:02:37"
$ws.Range("B316").Value = 0.0034

$ws.Range("A317").Value = "2023-12-14 12:02:54"
$ws.Range("B317").Value = 0.001

$ws.Range("A318").Value = "2023-12-14 12:03:03"
$ws.Range("B318").Value = 0.0002
